$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so Excel does not coerce these
# dotted / decimal-look-alike strings into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.013.90"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "1.884.20"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "244.04"
$ws.Range("E5").Value = "  -2.39%  "

$ws.Range("D6").Value = "0.9980"

$ws.Range("D7").Value = "0.4957"
$ws.Range("E7").Value = "  -0.58%  "

$ws.Range("D8").Value = "44.41"
$ws.Range("E8").Value = "  -2.22%  "

$ws.Range("D9").Value = "0.2918"
$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("D10").Value = "0.06628"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").Value = "1.880.63"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").Value = "17.02"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").Value = "0.07206"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("D14").Value = "0.6669"
$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "86.03"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "4.858"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").Value = "30.003.22"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "0.000007846"
$ws.Range("E18").Value = "  +3.91%  "

$ws.Range("D19").Value = "0.9981"

$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").Value = "2.122.70"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").Value = "0.9981"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "4.775"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").Value = "9.179"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").Value = "5.604"
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("D26").Value = "149.17"
$ws.Range("E26").Value = "  +2.73%  "

$ws.Range("D27").Value = "136.91"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").Value = "16.80"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "1.912"
$ws.Range("E29").Value = "  -2.44%  "

$ws.Range("D30").Value = "1.382"
$ws.Range("E30").Value = "  -0.73%  "

$ws.Range("D31").Value = "4.193"
$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("D32").Value = "0.08675"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").Value = "3.968"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").Value = "0.04994"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("D35").Value = "1.108"
$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").Value = "0.7033"
$ws.Range("E36").Value = "  +1.48%  "

$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").Value = "2.211"
$ws.Range("E38").Value = "  -5.38%  "

$ws.Range("D39").Value = "2.694"
$ws.Range("E39").Value = "  -2.13%  "

$ws.Range("D40").Value = "0.9358"
$ws.Range("E40").Value = "  -2.77%  "

$ws.Range("D41").Value = "0.01642"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("D42").Value = "5.980"
$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "0.4198"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "101.78"
$ws.Range("E45").Value = "  -2.34%  "

$ws.Range("D46").Value = "7.583"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("D47").Value = "0.1264"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "0.05728"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("D49").Value = "32.47"
$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").Value = "8.273"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").Value = "0.3712"
$ws.Range("E51").Value = "  -0.78%  "
